$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 6
$ws.Range("H6").Value = 896
$ws.Range("I6").Value = 896
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 2688
$ws.Range("L6").Value = 0
$ws.Range("M6").Value = -2576

# Row 15
$ws.Range("H15").Value = 13515203
$ws.Range("I15").Value = 13515203
$ws.Range("J15").Value = 0
$ws.Range("K15").Value = 40545609
$ws.Range("L15").Value = 0
$ws.Range("M15").Value = -40545440

# Row 17
$ws.Range("H17").Value = 1418.5454
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 1418.5454
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 4255.6362
$ws.Range("N17").Value = -4591.6362

# Row 32
$ws.Range("H32").Value = 55557600
$ws.Range("I32").Value = 100000720
$ws.Range("J32").Value = 3700.5
$ws.Range("K32").Value = 100000720
$ws.Range("L32").Value = 3700.5
$ws.Range("M32").Value = -100000394
$ws.Range("N32").Value = -4352.5

# Row 125
$ws.Range("H125").Value = 3116.0625
$ws.Range("I125").Value = 1916
$ws.Range("J125").Value = 3287.5
$ws.Range("K125").Value = 17244
$ws.Range("L125").Value = 29587.5
$ws.Range("M125").Value = -14784
$ws.Range("N125").Value = -34507.5

# Row 137
$ws.Range("H137").Value = 16353420
$ws.Range("I137").Value = 835209.5600000001
$ws.Range("J137").Value = 37044370
$ws.Range("K137").Value = 2505628.68
$ws.Range("L137").Value = 111133110
$ws.Range("M137").Value = -2503078.68
$ws.Range("N137").Value = -111138210

# Row 138
$ws.Range("H138").Value = 5437.2393
$ws.Range("I138").Value = 2386.4
$ws.Range("J138").Value = 5809.2925
$ws.Range("K138").Value = 7159.200000000001
$ws.Range("L138").Value = 17427.8775
$ws.Range("M138").Value = -2019.200000000001
$ws.Range("N138").Value = -27707.8775

$ws = $wb.Worksheets.Item("ARM")
# Row 61
$ws.Range("H61").Value = 38627.25
$ws.Range("I61").Value = 50836.332
$ws.Range("J61").Value = 2000
$ws.Range("K61").Value = 50836.332
$ws.Range("L61").Value = 2000
$ws.Range("M61").Value = -50624.332
$ws.Range("N61").Value = -2424

# Row 74
$ws.Range("H74").Value = 13159016
$ws.Range("I74").Value = 35715270
$ws.Range("J74").Value = 1199
$ws.Range("K74").Value = 35715270
$ws.Range("L74").Value = 1199
$ws.Range("M74").Value = -35714396
$ws.Range("N74").Value = -2947

# Row 77
$ws.Range("H77").Value = 13159016
$ws.Range("I77").Value = 35715270
$ws.Range("J77").Value = 1199
$ws.Range("K77").Value = 178576350
$ws.Range("L77").Value = 5995
$ws.Range("M77").Value = -178571982
$ws.Range("N77").Value = -14731

# Row 122
$ws.Range("H122").Value = 3939.25
$ws.Range("I122").Value = 2591
$ws.Range("J122").Value = 6785.5557
$ws.Range("K122").Value = 7773
$ws.Range("L122").Value = 20356.6671
$ws.Range("M122").Value = -5323
$ws.Range("N122").Value = -25256.6671

# Row 132
$ws.Range("H132").Value = 18130.715
$ws.Range("I132").Value = 21706.285
$ws.Range("J132").Value = 3828.4285
$ws.Range("K132").Value = 65118.855
$ws.Range("L132").Value = 11485.2855
$ws.Range("M132").Value = -62588.855
$ws.Range("N132").Value = -16545.2855

# Row 136
$ws.Range("H136").Value = 38627.25
$ws.Range("I136").Value = 50836.332
$ws.Range("J136").Value = 2000
$ws.Range("K136").Value = 152508.996
$ws.Range("L136").Value = 6000
$ws.Range("M136").Value = -149958.996
$ws.Range("N136").Value = -11100

$ws = $wb.Worksheets.Item("BSM")
# Row 86
$ws.Range("H86").Value = 1705.0714
$ws.Range("I86").Value = 1859
$ws.Range("J86").Value = 1499.8334
$ws.Range("K86").Value = 1859
$ws.Range("L86").Value = 1499.8334
$ws.Range("M86").Value = -736
$ws.Range("N86").Value = -3745.8334

# Row 89
$ws.Range("H89").Value = 1705.0714
$ws.Range("I89").Value = 1859
$ws.Range("J89").Value = 1499.8334
$ws.Range("K89").Value = 9295
$ws.Range("L89").Value = 7499.166999999999
$ws.Range("M89").Value = -3679
$ws.Range("N89").Value = -18731.167

# Row 105
$ws.Range("H105").Value = 1639.08
$ws.Range("I105").Value = 1360.15
$ws.Range("J105").Value = 2754.8
$ws.Range("K105").Value = 1360.15
$ws.Range("L105").Value = 2754.8
$ws.Range("M105").Value = 386.8499999999999
$ws.Range("N105").Value = -6248.8

$ws = $wb.Worksheets.Item("CRP")
# Row 4
$ws.Range("H4").Value = 5000
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 5000
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 5000
$ws.Range("N4").Value = -5224

# Row 31
$ws.Range("H31").Value = 5736.639
$ws.Range("I31").Value = 2607.5
$ws.Range("J31").Value = 6668.7236
$ws.Range("K31").Value = 2607.5
$ws.Range("L31").Value = 6668.7236
$ws.Range("M31").Value = -2312.5
$ws.Range("N31").Value = -7258.7236

# Row 34
$ws.Range("H34").Value = 5736.639
$ws.Range("I34").Value = 2607.5
$ws.Range("J34").Value = 6668.7236
$ws.Range("K34").Value = 2607.5
$ws.Range("L34").Value = 6668.7236
$ws.Range("M34").Value = -2405.5
$ws.Range("N34").Value = -7072.7236

# Row 99
$ws.Range("H99").Value = 16346
$ws.Range("I99").Value = 27004.4
$ws.Range("J99").Value = 10424.667
$ws.Range("K99").Value = 27004.4
$ws.Range("L99").Value = 10424.667
$ws.Range("M99").Value = -25506.4
$ws.Range("N99").Value = -13420.667

# Row 122
$ws.Range("H122").Value = 4129.174
$ws.Range("I122").Value = 2810.1667
$ws.Range("J122").Value = 5568.091
$ws.Range("K122").Value = 8430.500100000001
$ws.Range("L122").Value = 16704.273
$ws.Range("M122").Value = -5980.500100000001
$ws.Range("N122").Value = -21604.273

# Row 126
$ws.Range("H126").Value = 16346
$ws.Range("I126").Value = 27004.4
$ws.Range("J126").Value = 10424.667
$ws.Range("K126").Value = 81013.20000000001
$ws.Range("L126").Value = 31274.001
$ws.Range("M126").Value = -78543.20000000001
$ws.Range("N126").Value = -36214.001

$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value = 433051.38
$ws.Range("I5").Value = 582.6667
$ws.Range("J5").Value = 519545.12
$ws.Range("K5").Value = 1748.0001
$ws.Range("L5").Value = 1558635.36
$ws.Range("M5").Value = -1636.0001
$ws.Range("N5").Value = -1558859.36

# Row 29
$ws.Range("H29").Value = 150
$ws.Range("I29").Value = 0
$ws.Range("J29").Value = 150
$ws.Range("K29").Value = 0
$ws.Range("L29").Value = 450
$ws.Range("N29").Value = -1004

# Row 39
$ws.Range("H39").Value = 3235
$ws.Range("I39").Value = 4945
$ws.Range("J39").Value = 2665
$ws.Range("K39").Value = 14835
$ws.Range("L39").Value = 7995
$ws.Range("M39").Value = -14541
$ws.Range("N39").Value = -8583

# Row 98
$ws.Range("H98").Value = 918.4666999999999
$ws.Range("I98").Value = 225
$ws.Range("J98").Value = 968
$ws.Range("K98").Value = 675
$ws.Range("L98").Value = 2904
$ws.Range("M98").Value = 823
$ws.Range("N98").Value = -5900

# Row 132
$ws.Range("H132").Value = 8231.429
$ws.Range("I132").Value = 1030
$ws.Range("J132").Value = 17833.334
$ws.Range("K132").Value = 9270
$ws.Range("L132").Value = 160500.006
$ws.Range("M132").Value = -6740
$ws.Range("N132").Value = -165560.006

# Row 135
$ws.Range("H135").Value = 433051.38
$ws.Range("I135").Value = 582.6667
$ws.Range("J135").Value = 519545.12
$ws.Range("K135").Value = 5244.0003
$ws.Range("L135").Value = 4675906.08
$ws.Range("M135").Value = -2709.0003
$ws.Range("N135").Value = -4680976.08

# Row 137
$ws.Range("H137").Value = 6680524
$ws.Range("I137").Value = 1280
$ws.Range("J137").Value = 7708100
$ws.Range("K137").Value = 3840
$ws.Range("L137").Value = 23124300
$ws.Range("M137").Value = 1260
$ws.Range("N137").Value = -23134500

# Row 141
$ws.Range("H141").Value = 12477.9375
$ws.Range("I141").Value = 6465.2
$ws.Range("J141").Value = 22499.166
$ws.Range("K141").Value = 19395.6
$ws.Range("L141").Value = 67497.49800000001
$ws.Range("M141").Value = -14215.6
$ws.Range("N141").Value = -77857.49800000001

$ws = $wb.Worksheets.Item("GSM")
# Row 107
$ws.Range("H107").Value = 5953192
$ws.Range("I107").Value = 15873517
$ws.Range("J107").Value = 996.8
$ws.Range("K107").Value = 15873517
$ws.Range("L107").Value = 996.8
$ws.Range("M107").Value = -15871597
$ws.Range("N107").Value = -4836.8

# Row 132
$ws.Range("H132").Value = 416695
$ws.Range("I132").Value = 113767.78
$ws.Range("J132").Value = 912394.0600000001
$ws.Range("K132").Value = 341303.34
$ws.Range("L132").Value = 2737182.18
$ws.Range("M132").Value = -338773.34
$ws.Range("N132").Value = -2742242.18

$ws = $wb.Worksheets.Item("LTW")
# Row 16
$ws.Range("H16").Value = 6668403.5
$ws.Range("I16").Value = 10001320
$ws.Range("J16").Value = 2569.2
$ws.Range("K16").Value = 10001320
$ws.Range("L16").Value = 2569.2
$ws.Range("M16").Value = -10001150
$ws.Range("N16").Value = -2909.2

# Row 40
$ws.Range("H40").Value = 14494921
$ws.Range("I40").Value = 2233
$ws.Range("J40").Value = 41668708
$ws.Range("K40").Value = 2233
$ws.Range("L40").Value = 41668708
$ws.Range("M40").Value = -2097
$ws.Range("N40").Value = -41668980

# Row 131
$ws.Range("H131").Value = 55767.6
$ws.Range("I131").Value = 0
$ws.Range("J131").Value = 55767.6
$ws.Range("K131").Value = 0
$ws.Range("L131").Value = 55767.6
$ws.Range("N131").Value = -65847.60000000001

# Row 132
$ws.Range("H132").Value = 3344.2593
$ws.Range("I132").Value = 3396
$ws.Range("J132").Value = 2930.3333
$ws.Range("K132").Value = 10188
$ws.Range("L132").Value = 8790.999899999999
$ws.Range("M132").Value = -7658
$ws.Range("N132").Value = -13850.9999

# Row 136
$ws.Range("H136").Value = 2857.1765
$ws.Range("I136").Value = 2410.75
$ws.Range("J136").Value = 10000
$ws.Range("K136").Value = 7232.25
$ws.Range("L136").Value = 30000
$ws.Range("M136").Value = -4682.25
$ws.Range("N136").Value = -35100

$ws = $wb.Worksheets.Item("WVR")
# Row 6
$ws.Range("H6").Value = 0
$ws.Range("I6").Value = 0
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 0
$ws.Range("L6").Value = 0
$ws.Range("M6").ClearContents()

# Row 41
$ws.Range("H41").Value = 35000
$ws.Range("I41").Value = 0
$ws.Range("J41").Value = 35000
$ws.Range("K41").Value = 0
$ws.Range("L41").Value = 35000
$ws.Range("N41").Value = -35780

# Row 107
$ws.Range("H107").Value = 2499.7932
$ws.Range("I107").Value = 2492.8
$ws.Range("J107").Value = 2515.3333
$ws.Range("K107").Value = 7478.400000000001
$ws.Range("L107").Value = 7545.999899999999
$ws.Range("M107").Value = -5558.400000000001
$ws.Range("N107").Value = -11385.9999

# Row 122
$ws.Range("H122").Value = 3150.862
$ws.Range("I122").Value = 3373.5386
$ws.Range("J122").Value = 1221
$ws.Range("K122").Value = 10120.6158
$ws.Range("L122").Value = 3663
$ws.Range("M122").Value = -7670.6158
$ws.Range("N122").Value = -8563
